# testData.xlsx — "RegDetails" sheet gets a fresh row of sample
# credentials and loses the extra (3rd) data row, shrinking the shared
# strings table accordingly and tightening column B's width.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegDetails")

# Remove the 3rd row (A3/B3 - "testdata56"/"admin96") entirely, shifting
# the used range back down to A1:B2.
$ws.Range("A3:B3").EntireRow.Delete()

# Replace the remaining data row's values with the new sample data.
$ws.Range("A2").Value = "data45"
$ws.Range("B2").Value = "bread88"

# Column B gets an explicit width (previously auto/default).
$ws.Columns.Item(2).ColumnWidth = 9.63

# Move the selection to follow the data (previously B3, now gone).
$ws.Range("B2").Select()
